# Append two new daily rows ("2026-01-01" and "2026-01-02") to the main
# GSC export "Chart" sheet, each with 0 Non-HTTPS URLs and 29 HTTPS URLs,
# mirroring the existing rows that make up the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New data to append, in order.
$newDates = @("2026-01-01", "2026-01-02")
$httpsCounts = @(29, 29)

# Last currently populated row in the sheet (row 87 -> date 2025-12-31).
$lastRow = $ws.UsedRange.Rows.Count

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $lastRow + 1 + $i
    $dateCell = $ws.Cells.Item($row, 1)

    # A leading apostrophe forces the text to be stored as a literal string
    # instead of being auto-parsed into a date serial number, matching the
    # plain-text date cells already used throughout column A.
    $dateCell.Value = "'" + $newDates[$i]

    # Re-stamp the cell formatting from an existing data cell so the new
    # cell uses the same (default) style as the rest of the column rather
    # than the "quote prefix" style implicitly created above.
    $ws.Cells.Item(2, 1).Copy()
    $dateCell.PasteSpecial(-4122)

    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = $httpsCounts[$i]
}
